# Auto-generated: update cached market-price / profit figures across 8 job sheets.
# Source data refreshed by the scheduled market-data runner (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC: Morning Glass of Ether / Ether
$ws.Range("H15").Value = 16667975
$ws.Range("I15").Value = 16667975
$ws.Range("K15").Value = 50003925
$ws.Range("M15").Value = -50003756

# ALC: You Grow, Girl / Growth Formula Beta
$ws.Range("H18").Value = 906.2857
$ws.Range("I18").Value = 906.2857
$ws.Range("K18").Value = 906.2857
$ws.Range("M18").Value = -622.2857

# ALC: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 5491.857
$ws.Range("I40").Value = 1700
$ws.Range("J40").Value = 9283.714
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 9283.714
$ws.Range("M40").Value = -1525
$ws.Range("N40").Value = -9633.714

# ALC: Tomes Roam on the Range / Dhalmelskin Codex
$ws.Range("H75").Value = 66166.5
$ws.Range("J75").Value = 66166.5
$ws.Range("L75").Value = 66166.5
$ws.Range("N75").Value = -68038.5

# ALC: Field Trip to the Unknown (L) / Dhalmelskin Codex
$ws.Range("H78").Value = 66166.5
$ws.Range("J78").Value = 66166.5
$ws.Range("L78").Value = 198499.5
$ws.Range("N78").Value = -207859.5

# ALC: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 903384.8
$ws.Range("I132").Value = 903384.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2710154.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2707624.4
$ws.Range("N132").ClearContents()

# ALC: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1563
$ws.Range("I138").Value = 935.04346
$ws.Range("J138").Value = 2323.158
$ws.Range("K138").Value = 2805.13038
$ws.Range("L138").Value = 6969.474
$ws.Range("M138").Value = 2334.86962
$ws.Range("N138").Value = -17249.474

# ALC: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 70437.25
$ws.Range("J140").Value = 70437.25
$ws.Range("L140").Value = 70437.25
$ws.Range("N140").Value = -80797.25

$ws = $wb.Worksheets.Item("ARM")
# ARM: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3923637.5
$ws.Range("I132").Value = 1523.3846
$ws.Range("K132").Value = 4570.1538
$ws.Range("M132").Value = -2040.1538

$ws = $wb.Worksheets.Item("CRP")
# CRP: The Arsenal of Theocracy / Cobalt Halberd
$ws.Range("H50").Value = 39615.332
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 39615.332
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 39615.332
$ws.Range("N50").Value = -40865.332
$ws.Range("M50").ClearContents()

# CRP: Bowing to Greater Power / Yew Longbow
$ws.Range("H60").Value = 39938.824
$ws.Range("J60").Value = 42310
$ws.Range("L60").Value = 42310
$ws.Range("N60").Value = -43332

# CRP: License to Heal / Dark Chestnut Rod
$ws.Range("H74").Value = 96219.55499999999
$ws.Range("J74").Value = 96219.55499999999
$ws.Range("L74").Value = 96219.55499999999
$ws.Range("N74").Value = -97967.55499999999

# CRP: Purified Polyrhythm (L) / Dark Chestnut Rod
$ws.Range("H77").Value = 96219.55499999999
$ws.Range("J77").Value = 96219.55499999999
$ws.Range("L77").Value = 288658.665
$ws.Range("N77").Value = -297394.665

# CRP: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 11388.5
$ws.Range("I86").Value = 7733
$ws.Range("K86").Value = 7733
$ws.Range("M86").Value = -6610

# CRP: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 11388.5
$ws.Range("I89").Value = 7733
$ws.Range("K89").Value = 38665
$ws.Range("M89").Value = -33049

# CRP: Built to Last / White Oak Lumber
$ws.Range("H107").Value = 1463
$ws.Range("I107").Value = 1267.1052
$ws.Range("J107").Value = 1928.25
$ws.Range("K107").Value = 1267.1052
$ws.Range("L107").Value = 1928.25
$ws.Range("M107").Value = 652.8948
$ws.Range("N107").Value = -5768.25

# CRP: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 7043.625
$ws.Range("I132").Value = 5907
$ws.Range("K132").Value = 17721
$ws.Range("M132").Value = -15191

# CRP: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 10104462
$ws.Range("I134").Value = 11114458
$ws.Range("K134").Value = 33343374
$ws.Range("M134").Value = -33340839

$ws = $wb.Worksheets.Item("CUL")
# CUL: The Perks of Life at Sea / Jerked Beef
$ws.Range("H51").Value = 12802.5
$ws.Range("I51").Value = 6800
$ws.Range("J51").Value = 14003
$ws.Range("K51").Value = 20400
$ws.Range("L51").Value = 42009
$ws.Range("M51").Value = -19940
$ws.Range("N51").Value = -42929

# CUL: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 804.8333
$ws.Range("I107").Value = 715.25
$ws.Range("J107").Value = 849.625
$ws.Range("K107").Value = 2145.75
$ws.Range("L107").Value = 2548.875
$ws.Range("M107").Value = -225.75
$ws.Range("N107").Value = -6388.875

# CUL: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 929
$ws.Range("J122").Value = 932.625
$ws.Range("L122").Value = 8393.625
$ws.Range("N122").Value = -13293.625

# CUL: A Stickler for Carrots / Carrot Nibbles
$ws.Range("H127").Value = 50329.445
$ws.Range("J127").Value = 56433.125
$ws.Range("L127").Value = 169299.375
$ws.Range("N127").Value = -179219.375

$ws = $wb.Worksheets.Item("GSM")
# GSM: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 2587
$ws.Range("J2").Value = 3399.3333
$ws.Range("L2").Value = 3399.3333
$ws.Range("N2").Value = -3625.3333

# GSM: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 674866.7
$ws.Range("I122").Value = 674866.7
$ws.Range("K122").Value = 2024600.1
$ws.Range("M122").Value = -2022150.1

# GSM: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 954746.8
$ws.Range("I132").Value = 1251970.8
$ws.Range("K132").Value = 3755912.4
$ws.Range("M132").Value = -3753382.4

$ws = $wb.Worksheets.Item("LTW")
# LTW: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 4542.7
$ws.Range("I22").Value = 1301
$ws.Range("J22").Value = 6703.8335
$ws.Range("K22").Value = 1301
$ws.Range("L22").Value = 6703.8335
$ws.Range("M22").Value = -1006
$ws.Range("N22").Value = -7293.8335

# LTW: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 4542.7
$ws.Range("I27").Value = 1301
$ws.Range("J27").Value = 6703.8335
$ws.Range("K27").Value = 1301
$ws.Range("L27").Value = 6703.8335
$ws.Range("M27").Value = -1194
$ws.Range("N27").Value = -6917.8335

# LTW: Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 23812876
$ws.Range("I93").Value = 30306152
$ws.Range("J93").Value = 4201.3335
$ws.Range("K93").Value = 30306152
$ws.Range("L93").Value = 4201.3335
$ws.Range("M93").Value = -30304904
$ws.Range("N93").Value = -6697.3335

# LTW: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 10432.737
$ws.Range("I132").Value = 9077.375
$ws.Range("J132").Value = 17661.334
$ws.Range("K132").Value = 27232.125
$ws.Range("L132").Value = 52984.00199999999
$ws.Range("M132").Value = -24702.125
$ws.Range("N132").Value = -58044.00199999999

$ws = $wb.Worksheets.Item("WVR")
# WVR: After the Smock-down / Linen Smock
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# WVR: Where the Dragonflies, the Net Catches / Crawler Silk
$ws.Range("H81").Value = 3402603
$ws.Range("I81").Value = 3969504
$ws.Range("J81").Value = 1198
$ws.Range("K81").Value = 7939008
$ws.Range("L81").Value = 2396
$ws.Range("M81").Value = -7937947
$ws.Range("N81").Value = -4518

# WVR: To Kill a Dragon on Nameday (L) / Crawler Silk
$ws.Range("H84").Value = 3402603
$ws.Range("I84").Value = 3969504
$ws.Range("J84").Value = 1198
$ws.Range("K84").Value = 39695040
$ws.Range("L84").Value = 11980
$ws.Range("M84").Value = -39689736
$ws.Range("N84").Value = -22588

# WVR: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 1582.7307
$ws.Range("I113").Value = 1624.5714
$ws.Range("J113").Value = 1407
$ws.Range("K113").Value = 4873.7142
$ws.Range("L113").Value = 4221
$ws.Range("M113").Value = -2703.7142
$ws.Range("N113").Value = -8561

# WVR: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 3313.3333
$ws.Range("I126").Value = 3220.625
$ws.Range("J126").Value = 3498.75
$ws.Range("K126").Value = 9661.875
$ws.Range("L126").Value = 10496.25
$ws.Range("M126").Value = -7191.875
$ws.Range("N126").Value = -15436.25

# WVR: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1888.619
$ws.Range("I132").Value = 1763.6471
$ws.Range("K132").Value = 5290.9413
$ws.Range("M132").Value = -2760.9413

# WVR: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1628.579
$ws.Range("I136").Value = 1534.0769
$ws.Range("K136").Value = 4602.2307
$ws.Range("M136").Value = -2052.2307
